$wb = $excel.ActiveWorkbook

# --- 1. Reorder sheet tabs: review_info first, hotel_info second ---
$hotelInfo = $wb.Worksheets.Item("hotel_info")
$reviewInfo = $wb.Worksheets.Item("review_info")
$hotelInfo.Move($null, $reviewInfo)

# --- 2. Insert a new "State" column into hotel_info (between Hotel_Name and City) ---
$ws = $wb.Worksheets.Item("hotel_info")
$ws.Range("C1:C2").EntireColumn.Insert()
$ws.Cells.Item(1, 3).Value = "State"
$ws.Cells.Item(2, 3).Value = "Louisiana"
